$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("results")

# Insert a new row at position 2, shifting existing rows 2-8 down to 3-9.
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new metric.
$ws.Range("A2").Value = "id_credito_testeado"
$ws.Range("B2").Value = 169

# The inserted row doesn't automatically carry the label / value column
# formatting from the surrounding rows, so copy it explicitly from row 3
# (the row that used to be row 2 before the insert).
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# Update the big JSON-ish summary cell (now on row 9) with the richer
# per-entity breakdown.
$ws.Range("B9").Value = '{''ner-general'': {''entidades_correctas'': 0, ''entidades_incorrectas'': 0, ''porcentaje_entidades_correctas'': 1, ''paginas_correctas'': 0, ''paginas_incorrectas'': 0, ''porcentaje_paginas_correctas'': 1}, ''ner-escritura'': {''entidades_correctas'': 34, ''entidades_incorrectas'': 3, ''porcentaje_entidades_correctas'': 0.918918918918919, ''paginas_correctas'': 34, ''paginas_incorrectas'': 3, ''porcentaje_paginas_correctas'': 0.918918918918919, ''CIUDAD_INMUEBLE'': {''pagina_encontrado'': 27, ''pagina_es_correcta'': True, ''pagina_real'': 27, ''texto_es_correcto'': True, ''valor_encontrado'': ''CHIHUAHUA'', ''valor_real'': ''CHIHUAHUA''}, ''CLAUSULA_ACTUALIZACION_CREDITO'': {''pagina_encontrado'': 10, ''pagina_es_correcta'': True, ''pagina_real'': 10, ''texto_es_correcto'': True, ''valor_encontrado'': ''PRIMERA'', ''valor_real'': ''PRIMERA''}, ''CLAUSULA_CONSTITUCION_HIPOTECA'': {''pagina_encontrado'': 19, ''pagina_es_correcta'': True, ''pagina_real'': 19, ''texto_es_correcto'': True, ''valor_encontrado'': ''VIGESIMA TERCERA'', ''valor_real'': ''VIGESIMA TERCERA''}, ''CLAUSULA_JURISDICCION_COMPETENCIA'': {''pagina_encontrado'': 20, ''pagina_es_correcta'': True, ''pagina_real'': 20, ''texto_es_correcto'': True, ''valor_encontrado'': ''VIGESIMA QUINTA'', ''valor_real'': ''VIGESIMA QUINTA''}, ''CLAUSULA_PRORROGA'': {''pagina_encontrado'': 14, ''pagina_es_correcta'': True, ''pagina_real'': 14, ''texto_es_correcto'': True, ''valor_encontrado'': ''DECIMA SEGUNDA'', ''valor_real'': ''DECIMA SEGUNDA''}, ''CLAUSULA_VENCIMIENTO_ANTICIPADO'': {''pagina_encontrado'': 17, ''pagina_es_correcta'': True, ''pagina_real'': 17, ''texto_es_correcto'': True, ''valor_encontrado'': ''VIGESIMA'', ''valor_real'': ''VIGESIMA''}, ''CLAUSULA_INT_ORD'': {''pagina_encontrado'': 13, ''pagina_es_correcta'': True, ''pagina_real'': 13, ''texto_es_correcto'': True, ''valor_encontrado'': ''NOVENA'', ''valor_real'': ''NOVENA''}, ''CLAUSULA_PLAZO'': {''pagina_encontrado'': 11, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''SEXTA'', ''valor_real'': ''SEXTA''}, ''CREDITO_PESOS'': {''pagina_encontrado'': 10, ''pagina_es_correcta'': True, ''pagina_real'': 10, ''texto_es_correcto'': True, ''valor_encontrado'': ''$244,999.92 (DOSCIENTOS CUARENTA Y CUATRO MIL NOVECIENTOS NOVENTA Y NUEVE PESOS 92/100)'', ''valor_real'': ''$244,999.92 (DOSCIENTOS CUARENTA Y CUATRO MIL NOVECIENTOS NOVENTA Y NUEVE PESOS 92/100)''}, ''DESCRIPCION_INMUEBLE'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''LOTE DE TERRENO NUMERO 3 TRES, DE LA MANZANA 9 NUEVE'', ''valor_real'': ''LOTE DE TERRENO NUMERO 3 TRES, DE LA MANZANA 9 NUEVE''}, ''DISTRITO_RPP'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': False, ''pagina_real'': 62, ''texto_es_correcto'': False, ''valor_encontrado'': ''CHIHUAHUA,'', ''valor_real'': ''CHIHUAHUA, CHIH''}, ''DOMICILIO_GARANTIA'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''FINCA MARCADA CON EL NUMERO 11,504 ONCE MIL QUINIENTOS CUATRO, DE LA CALLE DEL LUCERO, DEL FRACCIONAMIENTO LA GALERA ETAPA V CINCO, DE ESTA CIUDAD DE CHIHUAHUA, CONSTRUIDA SOBRE EL LOTE DE TERRENO NUMERO 3 TRES, DE LA MANZANA 9 NUEVE'', ''valor_real'': ''FINCA MARCADA CON EL NUMERO 11,504 ONCE MIL QUINIENTOS CUATRO, DE LA CALLE DEL LUCERO, DEL FRACCIONAMIENTO LA GALERA ETAPA V CINCO, DE ESTA CIUDAD DE CHIHUAHUA, CONSTRUIDA SOBRE EL LOTE DE TERRENO NUMERO 3 TRES, DE LA MANZANA 9 NUEVE''}, ''ESTADO_INMUEBLE'': {''pagina_encontrado'': 27, ''pagina_es_correcta'': True, ''pagina_real'': 27, ''texto_es_correcto'': True, ''valor_encontrado'': ''CHIHUAHUA'', ''valor_real'': ''CHIHUAHUA''}, ''FECHA_ESCRITURA'': {''pagina_encontrado'': 4, ''pagina_es_correcta'': True, ''pagina_real'': 4, ''texto_es_correcto'': True, ''valor_encontrado'': ''ONCE DIAS DEL MES DE JULIO DEL AÑO DOS MIL CINCO'', ''valor_real'': ''ONCE DIAS DEL MES DE JULIO DEL AÑO DOS MIL CINCO''}, ''JURISDICCION_NOTARIO'': {''pagina_encontrado'': 4, ''pagina_es_correcta'': True, ''pagina_real'': 4, ''texto_es_correcto'': True, ''valor_encontrado'': ''CHIHUAHUA, CHIHUAHUA'', ''valor_real'': ''CHIHUAHUA, CHIHUAHUA''}, ''JURISDICCION_SOMETIMIENTO_EXPRESO'': {''pagina_encontrado'': 20, ''pagina_es_correcta'': True, ''pagina_real'': 20, ''texto_es_correcto'': True, ''valor_encontrado'': ''DISTRITO FEDERAL'', ''valor_real'': ''DISTRITO FEDERAL''}, ''COLINDANCIAS_1'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''POR SU FRENTE 7.00 MTS. SIETE METROS CON LA CALLE DEL LUCERO'', ''valor_real'': ''POR SU FRENTE 7.00 MTS. SIETE METROS CON LA CALLE DEL LUCERO''}, ''COLINDANCIAS_2'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''POR SU COSTADO DERECHO 17.15 MTS. DIECISIETE METROS QUINCE CENTIMETROS CON LOTE NUMERO 4 CUATRO'', ''valor_real'': ''POR SU COSTADO DERECHO 17.15 MTS. DIECISIETE METROS QUINCE CENTÍMETROS CON LOTE NÚMERO 4 CUATRO''}, ''COLINDANCIAS_3'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''POR SU COSTADO IZQUIERDO 17.15 MTS. DIECISIETE METROS QUINCE CENTIMETROS CON LOTE NUMERO 2 DOS'', ''valor_real'': ''POR SU COSTADO IZQUIERDO 17.15 MTS. DIECISIETE METROS QUINCE CENTÍMETROS CON LOTE NÚMERO 2 DOS''}, ''COLINDANCIAS_4'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''POR SU ESPALDA 7.00 MTS. SIETE METROS CON LOTE NUMERO 40 CUARENTA'', ''valor_real'': ''POR SU ESPALDA 7.00 MTS. SIETE METROS CON LOTE NÚMERO 40 CUARENTA''}, ''NOTARIO'': {''pagina_encontrado'': 4, ''pagina_es_correcta'': True, ''pagina_real'': 4, ''texto_es_correcto'': True, ''valor_encontrado'': ''LUIS ARTURO CALDERON TRUEBA'', ''valor_real'': ''LUIS ARTURO CALDERON TRUEBA''}, ''NUM_ CLAUSULA_DESCUENTO_SALARIO_REA'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''SEPTIMA'', ''valor_real'': ''SEPTIMA''}, ''NUM_ CLAUSULA_DESCUENTO_SALARIO_ROA'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''SEPTIMA'', ''valor_real'': ''SEPTIMA''}, ''NUM_CLAUSULA_INTERES_MORATORIO'': {''pagina_encontrado'': 13, ''pagina_es_correcta'': False, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''DECIMA'', ''valor_real'': ''DECIMA''}, ''NUM_CLAUSULA_OTORGAMIENTO_CREDITO'': {''pagina_encontrado'': 10, ''pagina_es_correcta'': True, ''pagina_real'': 10, ''texto_es_correcto'': True, ''valor_encontrado'': ''PRIMERA'', ''valor_real'': ''PRIMERA''}, ''NUM_NOTARIO'': {''pagina_encontrado'': 4, ''pagina_es_correcta'': True, ''pagina_real'': 4, ''texto_es_correcto'': True, ''valor_encontrado'': ''DIEZ'', ''valor_real'': ''DIEZ''}, ''PLAZO_ANIOS'': {''pagina_encontrado'': 11, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''30 TREINTA AÑOS'', ''valor_real'': ''30 TREINTA AÑOS''}, ''REA_VSM'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 13, ''texto_es_correcto'': True, ''valor_encontrado'': ''50.6160 CINCUENTA PUNTO SEIS UNO SEIS CERO'', ''valor_real'': ''50.6160 CINCUENTA PUNTO SEIS UNO SEIS CERO''}, ''PORCENTAJE_ROA'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''20% (VEINTE POR CIENTO)'', ''valor_real'': ''20% (VEINTE POR CIENTO)''}, ''ROA_VSM'': {''pagina_encontrado'': 12, ''pagina_es_correcta'': True, ''pagina_real'': 12, ''texto_es_correcto'': True, ''valor_encontrado'': ''41.8760 CUARENTA Y UNO PUNTO OCHO SIETE SEIS CERO'', ''valor_real'': ''41.8760 CUARENTA Y UNO PUNTO OCHO SIETE SEIS CERO''}, ''SUPERFICIE_CONSTRUCCION'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': False, ''valor_encontrado'': ''68,711.89 SESENTA Y OCHO MIL SETECIENTOS ONCE METROS OCHENTA Y NUEVE DECIMETROS CUADRADOS'', ''valor_real'': ''120.05 M2. CIENTO VEINTE METROS CINCO DECIMETROS CUADRADOS''}, ''SUPERFICIE_TERRENO'': {''pagina_encontrado'': 5, ''pagina_es_correcta'': True, ''pagina_real'': 5, ''texto_es_correcto'': True, ''valor_encontrado'': ''120.05 M2. CIENTO VEINTE METROS CINCO DECIMETROS CUADRADOS'', ''valor_real'': ''120.05 M2. CIENTO VEINTE METROS CINCO DECIMETROS CUADRADOS''}, ''TASA_INTERES_MORATORIO'': {''pagina_encontrado'': 49, ''pagina_es_correcta'': False, ''pagina_real'': 14, ''texto_es_correcto'': False, ''valor_encontrado'': ''9.00% ANUAL FIJO'', ''valor_real'': ''9.00% NUEVE POR CIENTO''}, ''TASA_INTERES_ORDINARIO'': {''pagina_encontrado'': 13, ''pagina_es_correcta'': True, ''pagina_real'': 13, ''texto_es_correcto'': True, ''valor_encontrado'': ''9.00% NUEVE PUNTO CERO POR CIENTO'', ''valor_real'': ''9.00% NUEVE PUNTO CERO POR CIENTO''}, ''TEXTO_ACTUALIZACION'': {''pagina_encontrado'': 10, ''pagina_es_correcta'': True, ''pagina_real'': 10, ''texto_es_correcto'': True, ''valor_encontrado'': ''EL TRABAJADOR RECONOCE DEBER Y SE OBLIGA A PAGAR A EL INFONAVIT EL MONTO DEL CREDITO OTORGADO EN ESTE ACTO, EN LOS TERMINOS Y CONDICIONES PACTADAS EN ESTA ESCRITURA, ACEPTANDO QUE EL SALDO DEL CREDITO SE INCREMENTARA EN LA MISMA PROPORCION EN QUE AUMENTE EL SALARIO MINIMO DIARIO GENERAL QUE RIJA EN EL DISTRITO FEDERAL'', ''valor_real'': ''EL TRABAJADOR RECONOCE DEBER Y SE OBLIGA A PAGAR A EL INFONAVIT EL MONTO DEL CREDITO OTORGADO EN ESTE ACTO, EN LOS TERMINOS Y CONDICIONES PACTADAS EN ESTA ESCRITURA, ACEPTANDO QUE EL SALDO DEL CREDITO SE INCREMENTARA EN LA MISMA PROPORCION EN QUE AUMENTE EL SALARIO MINIMO DIARIO GENERAL QUE RIJA EN EL DISTRITO FEDERAL''}, ''TEXTO_CLAUSULA_HIPOTECA'': {''pagina_encontrado'': 19, ''pagina_es_correcta'': True, ''pagina_real'': 19, ''texto_es_correcto'': True, ''valor_encontrado'': ''EL TRABAJADOR Y SU CONYUGE, SENOR(A) MARIA TERESA BARRIOS IBARRA, PARA GARANTIZAR TODAS Y CADA UNA DE LAS OBLIGACIONES QUE SE CONTRAEN POR VIRTUD DEL PRESENTE CONTRATO, CONSTITUYE(N) HIPOTECA, SIENDO LAS 10:30 DIEZ HORAS TREINTA MINUTOS DEL DIA DE LA FECHA, A FAVOR DE EL INFONAVIT SOBRE LA FINCA MARCADA CON EL NUMERO 11,504 ONCE MIL QUINIENTOS CUATRO, DE LA CALLE DEL LUCERO, DEL FRACCIONAMIENTO LA GALERA ETAPA V CINCO, DE ESTA CIUDAD DE CHIHUAHUA, CONSTRUIDA SOBRE EL LOTE DE TERRENO NUMERO 3 TRES, DE LA MANZANA 9 NUEVE, CON SUPERFICIE DE 120.05 M2. CIENTO VEINTE METROS CINCO DECIMETROS CUADRADOS, INMUEBLE RELACIONADO EN EL ANTECEDENTE PRIMERO DE ESTA ESCRITURA'', ''valor_real'': ''EL TRABAJADOR Y SU CONYUGE, SENOR(A) MARIA TERESA BARRIOS IBARRA, PARA GARANTIZAR TODAS Y CADA UNA DE LAS OBLIGACIONES QUE SE CONTRAEN POR VIRTUD DEL PRESENTE CONTRATO, CONSTITUYE(N) HIPOTECA, SIENDO LAS 10:30 DIEZ HORAS TREINTA MINUTOS DEL DIA DE LA FECHA, A FAVOR DE EL INFONAVIT SOBRE LA FINCA MARCADA CON EL NUMERO 11,504 ONCE MIL QUINIENTOS CUATRO, DE LA CALLE DEL LUCERO, DEL FRACCIONAMIENTO LA GALERA ETAPA V CINCO, DE ESTA CIUDAD DE CHIHUAHUA, CONSTRUIDA SOBRE EL LOTE DE TERRENO NUMERO 3 TRES, DE LA MANZANA 9 NUEVE, CON SUPERFICIE DE 120.05 M2. CIENTO VEINTE METROS CINCO DECIMETROS CUADRADOS, INMUEBLE RELACIONADO EN EL ANTECEDENTE PRIMERO DE ESTA ESCRITURA''}, ''TIPO_INMUEBLE'': {''pagina_encontrado'': 27, ''pagina_es_correcta'': True, ''pagina_real'': 27, ''texto_es_correcto'': True, ''valor_encontrado'': ''CASA HABITACION'', ''valor_real'': ''CASA HABITACION''}}, ''ner-certificado'': {''entidades_correctas'': 7, ''entidades_incorrectas'': 4, ''porcentaje_entidades_correctas'': 0.6363636363636364, ''paginas_correctas'': 7, ''paginas_incorrectas'': 4, ''porcentaje_paginas_correctas'': 1.0, ''CONDICIONES _CREDITO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''VSM'', ''valor_real'': ''VSM''}, ''CURP_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''OIRF790826HCHLML''}, ''DELEGACION'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''CHIHUAHUA'', ''valor_real'': ''CHIHUAHUA''}, ''FECHA_EMISION'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''25 DE AGOSTO DE 2017 ''}, ''MESES_ADEUDADOS_CERTIFICADO_ADEUDOS'': {''pagina_encontrado'': 9, ''pagina_es_correcta'': True, ''pagina_real'': 9, ''texto_es_correcto'': True, ''valor_encontrado'': ''5'', ''valor_real'': ''5''}, ''MONTO_VSM_INTERES_MORATORIO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''5.2000'', ''valor_real'': ''5.2000''}, ''NOMBRE_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''OLIVAS RAMIREZ FELIPE DE JESUS'', ''valor_real'': ''OLIVAS RAMIREZ FELIPE DE JESUS''}, ''NUMERO_CREDITO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': '' '', ''valor_real'': ''805117974''}, ''NUMERO_HOJAS_CA'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''9'', ''valor_real'': ''9''}, ''RFC_ACREDITADO'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': False, ''valor_encontrado'': ''FACULTADES'', ''valor_real'': ''OIRF790828152''}, ''SUERTE_PRINCIPAL'': {''pagina_encontrado'': 1, ''pagina_es_correcta'': True, ''pagina_real'': 1, ''texto_es_correcto'': True, ''valor_encontrado'': ''114.9030'', ''valor_real'': ''114.9030''}}, ''avaluo_adjudicacion'': {''entidades_correctas'': 0, ''entidades_incorrectas'': 0, ''porcentaje_entidades_correctas'': 1, ''paginas_correctas'': 0, ''paginas_incorrectas'': 0, ''porcentaje_paginas_correctas'': 1}}'
